$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 5900
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5900
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5900
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -6238
$ws.Range("H18").Value = 4917.6
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H43").Value = 1400
$ws.Range("I43").Value = 1400
$ws.Range("K43").Value = 1400
$ws.Range("M43").Value = -1331
$ws.Range("H53").Value = 375.5
$ws.Range("I53").Value = 291.8
$ws.Range("J53").Value = 459.2
$ws.Range("K53").Value = 291.8
$ws.Range("L53").Value = 459.2
$ws.Range("M53").Value = 345.2
$ws.Range("N53").Value = -1733.2
$ws.Range("H125").Value = 15000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 18000
$ws.Range("M125").Value = -15540
$ws.Range("H132").Value = 2332.889
$ws.Range("I132").Value = 2332.889
$ws.Range("K132").Value = 6998.667
$ws.Range("M132").Value = -4468.667
$ws.Range("H135").Value = 464.1111
$ws.Range("I135").Value = 464.1111
$ws.Range("K135").Value = 4176.9999
$ws.Range("M135").Value = -1641.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2945
$ws.Range("I2").Value = 2945
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2945
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2832
$ws.Range("N2").ClearContents()
$ws.Range("H4").Value = 899.3333
$ws.Range("I4").Value = 849.5
$ws.Range("J4").Value = 999
$ws.Range("K4").Value = 849.5
$ws.Range("L4").Value = 999
$ws.Range("M4").Value = -733.5
$ws.Range("N4").Value = -1231
$ws.Range("H5").Value = 434.14285
$ws.Range("J5").Value = 121
$ws.Range("L5").Value = 121
$ws.Range("N5").Value = -345
$ws.Range("H6").Value = 379500
$ws.Range("I6").Value = 379500
$ws.Range("K6").Value = 379500
$ws.Range("M6").Value = -379327
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H74").Value = 4349.7144
$ws.Range("I74").Value = 3574.6667
$ws.Range("K74").Value = 3574.6667
$ws.Range("M74").Value = -2700.6667
$ws.Range("H77").Value = 4349.7144
$ws.Range("I77").Value = 3574.6667
$ws.Range("K77").Value = 17873.3335
$ws.Range("M77").Value = -13505.3335
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H102").Value = 1473.3334
$ws.Range("I102").Value = 1473.3334
$ws.Range("K102").Value = 1473.3334
$ws.Range("M102").Value = 148.6666
$ws.Range("H116").Value = 2945
$ws.Range("I116").Value = 2945
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2945
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -651
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 2558.2222
$ws.Range("I132").Value = 1861.2858
$ws.Range("J132").Value = 4997.5
$ws.Range("K132").Value = 5583.857400000001
$ws.Range("L132").Value = 14992.5
$ws.Range("M132").Value = -3053.857400000001
$ws.Range("N132").Value = -20052.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2945
$ws.Range("I3").Value = 2945
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2945
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2831
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 434.14285
$ws.Range("J4").Value = 121
$ws.Range("L4").Value = 121
$ws.Range("N4").Value = -351
$ws.Range("H134").Value = 2190.7058
$ws.Range("I134").Value = 2149.3845
$ws.Range("K134").Value = 6448.1535
$ws.Range("M134").Value = -3913.1535

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2198.4
$ws.Range("I58").Value = 2000
$ws.Range("K58").Value = 2000
$ws.Range("M58").Value = -1797
$ws.Range("H132").Value = 3374.923
$ws.Range("J132").Value = 3065.6667
$ws.Range("L132").Value = 9197.000100000001
$ws.Range("N132").Value = -14257.0001
$ws.Range("H134").Value = 2077.7144
$ws.Range("I134").Value = 2158.3044
$ws.Range("J134").Value = 1707
$ws.Range("K134").Value = 6474.9132
$ws.Range("L134").Value = 5121
$ws.Range("M134").Value = -3939.9132
$ws.Range("N134").Value = -10191
$ws.Range("H136").Value = 2198.4
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450
$ws.Range("H141").Value = 44959.168
$ws.Range("J141").Value = 44959.168
$ws.Range("L141").Value = 44959.168
$ws.Range("N141").Value = -55319.168

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 31272108
$ws.Range("I11").Value = 62544068
$ws.Range("J11").Value = 146.66667
$ws.Range("K11").Value = 187632204
$ws.Range("L11").Value = 440.00001
$ws.Range("M11").Value = -187632064
$ws.Range("N11").Value = -720.00001
$ws.Range("H38").Value = 70
$ws.Range("I38").Value = 70
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 210
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 137
$ws.Range("N38").ClearContents()
$ws.Range("H56").Value = 11103.952
$ws.Range("I56").Value = 11103.952
$ws.Range("K56").Value = 11103.952
$ws.Range("M56").Value = -10573.952

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 49990
$ws.Range("J121").Value = 49990
$ws.Range("L121").Value = 49990
$ws.Range("N121").Value = -53484

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 8258.308000000001
$ws.Range("I122").Value = 7915
$ws.Range("K122").Value = 23745
$ws.Range("M122").Value = -21295
$ws.Range("H132").Value = 9497.875
$ws.Range("I132").Value = 12976.8
$ws.Range("K132").Value = 38930.39999999999
$ws.Range("M132").Value = -36400.39999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 22329.428
$ws.Range("J74").Value = 22291.834
$ws.Range("L74").Value = 22291.834
$ws.Range("N74").Value = -24163.834
$ws.Range("H77").Value = 22329.428
$ws.Range("J77").Value = 22291.834
$ws.Range("L77").Value = 66875.50199999999
$ws.Range("N77").Value = -76235.50199999999
$ws.Range("H122").Value = 1694.4
$ws.Range("I122").Value = 1421.4286
$ws.Range("K122").Value = 4264.2858
$ws.Range("M122").Value = -1814.2858
$ws.Range("H132").Value = 2107.6155
$ws.Range("I132").Value = 2500
$ws.Range("J132").Value = 1224.75
$ws.Range("K132").Value = 7500
$ws.Range("L132").Value = 3674.25
$ws.Range("M132").Value = -4970
$ws.Range("N132").Value = -8734.25
$ws.Range("H136").Value = 2605.375
$ws.Range("I136").Value = 2940.6667
$ws.Range("J136").Value = 1599.5
$ws.Range("K136").Value = 8822.000100000001
$ws.Range("L136").Value = 4798.5
$ws.Range("M136").Value = -6272.000100000001
$ws.Range("N136").Value = -9898.5
